# Adds a German translation column ("details_de") to the "vaccines links"
# sheet, renames the existing "details" column to "details_en", tweaks a
# couple of the English blurbs, and makes the "vaccines links" sheet the
# active tab/selection — per commit "added translation for vaccine info".

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("vaccines links")

# --- insert a new column C (pushes the old "link" column from C to D) ---
$ws2.Columns.Item(3).Insert() | Out-Null

# Column C is brand new text (German), give it roughly the same width as
# column B (details_en) rather than leaving it at the inherited "link" width.
$ws2.Columns.Item(3).ColumnWidth = 67.33

# --- headers ---
$ws2.Range("B1").Value = "details_en"
$ws2.Range("C1").Value = "details_de"

# --- row 2: BAG ---
$ws2.Range("A2").Value = "Bundesamt für Statistik (BAG)"
$ws2.Range("B2").Value = "information about vaccines and agreements with vaccine producers"
$ws2.Range("C2").Value = "Information zu den Impfstoffen und Verträgen mit Herstellern"

# --- row 3: SRF ---
$ws2.Range("B3").Value = "information about the vaccines"
$ws2.Range("C3").Value = "Informationen zu den Impfstoffen"

# --- row 4: Admin.ch ---
$ws2.Range("C4").Value = "press release: Janssen-Cilag AG reicht Zulassungsgesuch für ihren Impfstoffkandidaten ein"

# --- row 5: swissmedic ---
$ws2.Range("C5").Value = "wie und warum Impfstoffe im Körper wirken"

# --- row 6: swissinfo.ch ---
$ws2.Range("C6").Value = "How vaccine technology, choice and supply work in Switzerland (english)"

# --- row 7: nature comment ---
$ws2.Range("C7").Value = "Next-generation vaccine platforms for COVID-19 (english)"

# --- move the hyperlink that used to sit on C2 onto the new D2 location ---
$bagUrl = "https://www.bag.admin.ch/bag/en/home/krankheiten/ausbrueche-epidemien-pandemien/aktuelle-ausbrueche-epidemien/novel-cov/impfen.html"
$ws2.Range("C2").Hyperlinks.Delete() | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), $bagUrl) | Out-Null

# --- make "vaccines links" the active/selected sheet, matching the saved view state ---
$ws2.Activate()
$ws2.Range("C10").Select() | Out-Null
